$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new blank row at the top of the table's data region (worksheet row 2),
# which shifts all existing data rows (and their formatting) down by one, exactly
# like using Excel's "Insert Table Rows Above" on the first data row.
$ws.Rows("2:2").Insert()
$lo.Resize($ws.Range("A1:E21"))

# The row insert leaves the calculated-column formula in the very last row in a
# broken, self-referential state - restore it to the normal structured reference.
$ws.Range("D21").Formula = "=VIC_Mystery_cases[[#This Row],[Date]]"
$ws.Range("E21").Formula = "=VIC_Mystery_cases[[#This Row],[Date]]+14"

# Populate the new row with the new mystery-case record (18 Oct 2020 update).
$ws.Range("A2").Value = 44119
$ws.Range("B2").Value = 3081
$ws.Range("C2").Value = "https://www.dhhs.vic.gov.au/coronavirus-update-for-victoria-sunday-18-october-2020"
$ws.Range("D2").Formula = "=VIC_Mystery_cases[[#This Row],[Date]]"
$ws.Range("E2").Formula = "=VIC_Mystery_cases[[#This Row],[Date]]+14"

# Give the new row the same formatting as the rest of the data rows (the insert
# defaulted it to the header row's style).
$ws.Range("A3:E3").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# The existing hyperlink stayed anchored on the old C2 cell; move it to the cell
# that now holds that news-link text (C3).
$ws.Range("C2").Hyperlinks.Delete()
$ws.Range("C3").Hyperlinks.Add($ws.Range("C3"), "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-17-october-2020") | Out-Null

# Adding the hyperlink re-applied Excel's built-in "Hyperlink" style to C3; restore
# the original (non-hyperlink) cell formatting used by the rest of the column.
$ws.Range("C4").Copy()
$ws.Range("C3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Match the selection shown after the edit.
$ws.Range("A2").Select() | Out-Null

$excel.CalculateFullRebuild() | Out-Null
